# Applies the edits described by the diff: splits several runs so
# mis-spelled / re-punctuated words get their own <w:proofErr>-wrapped
# run, adds <w:strike/> formatting in a couple of spots, and appends a
# new strike-through bullet ("Writing the mainloop and the clean data
# function") at the end of the list.

$d = $word.ActiveDocument

function Set-ParagraphXml($paragraphIndex, $xmlFragment) {
    # Replace paragraph $paragraphIndex (1-based) contents in-place with
    # $xmlFragment (a well-formed <w:p>...</w:p> string), preserving the
    # paragraph's identity/attributes supplied inside the fragment itself.
    $totalBefore = $d.Paragraphs.Count
    $target = $d.Paragraphs($paragraphIndex).Range
    $payload = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $xmlFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($payload)
    # When $paragraphIndex was the very last paragraph in the body, Word
    # leaves a stray empty paragraph (carrying the old paragraph mark)
    # behind; trim it back off so the paragraph count is unchanged.
    if ($paragraphIndex -eq $totalBefore -and $d.Paragraphs.Count -gt $totalBefore) {
        $stray = $d.Paragraphs($d.Paragraphs.Count)
        $d.Range($stray.Range.Start - 1, $stray.Range.End).Delete()
    }
}

Set-ParagraphXml 4 '<w:p w14:paraId="7C4C9DB9" w14:textId="77777777" w:rsidR="00155D62" w:rsidRDefault="00155D62" w:rsidP="00155D62"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>tmr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> first test what does the initial guess look like and then start to do the fit for the 0V individually first.</w:t></w:r></w:p>'
Set-ParagraphXml 6 '<w:p w14:paraId="602229C4" w14:textId="3216AB74" w:rsidR="00144544" w:rsidRDefault="00155D62" w:rsidP="00155D62"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Finished plotting the initial guess directly, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>more</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> or less looks fine.</w:t></w:r></w:p>'
Set-ParagraphXml 10 '<w:p w14:paraId="6D06D265" w14:textId="31067535" w:rsidR="00D54C91" w:rsidRDefault="00D54C91" w:rsidP="00155D62"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Change </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>globalFit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> initial guess from list to class and use class elements</w:t></w:r></w:p>'
Set-ParagraphXml 14 '<w:p w14:paraId="51DF734B" w14:textId="02A3937D" w:rsidR="007E4F0A" w:rsidRDefault="00045C6E" w:rsidP="007E4F0A"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Need to change the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mainloop</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to control how to use the initial guess algo and pass it into the global fit.</w:t></w:r></w:p>'
Set-ParagraphXml 21 '<w:p w14:paraId="05EE9783" w14:textId="777F1C91" w:rsidR="00626E28" w:rsidRDefault="00626E28" w:rsidP="00626E28"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r w:rsidRPr="00951A10"><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Test new model’s performance on the original </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00951A10"><w:rPr><w:strike/></w:rPr><w:t>fits( individual</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00951A10"><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> not 0</w:t></w:r><w:r w:rsidRPr="00951A10"><w:rPr><w:rFonts w:hint="eastAsia"/><w:strike/></w:rPr><w:t>V</w:t></w:r><w:r w:rsidRPr="00951A10"><w:rPr><w:strike/></w:rPr><w:t>)</w:t></w:r></w:p>'
Set-ParagraphXml 25 '<w:p w14:paraId="20C49FF1" w14:textId="7294AAB0" w:rsidR="007725A7" w:rsidRDefault="007725A7"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r w:rsidRPr="00A16BD1"><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Make independent fit do not depend on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00A16BD1"><w:rPr><w:strike/></w:rPr><w:t>Vbi</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
Set-ParagraphXml 26 '<w:p w14:paraId="4C650DC2" w14:textId="29F7C83E" w:rsidR="00055E2B" w:rsidRPr="00055E2B" w:rsidRDefault="00055E2B" w:rsidP="00055E2B"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Fix slider and fix value</w:t></w:r></w:p>'
Set-ParagraphXml 27 '<w:p w14:paraId="21C9BA02" w14:textId="017474CC" w:rsidR="00055E2B" w:rsidRPr="00A16BD1" w:rsidRDefault="00055E2B" w:rsidP="00055E2B"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Fix value use vary in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>lmfit</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

# Append a brand-new bullet after the last paragraph: "Writing the mainloop and the clean data function"
$lastIndex = $d.Paragraphs.Count
$d.Paragraphs($lastIndex).Range.InsertParagraphAfter()
Set-ParagraphXml ($lastIndex + 1) '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Writing the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>mainloop</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> and the clean data function</w:t></w:r></w:p>'

Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
